# Insert a new data row before current row 50 (shifting existing rows 50-85 down to 51-86)
# and populate the new row with the latest weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 50; existing rows shift down.
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new record.
$ws.Cells.Item(50, 1).Value = 5
$ws.Cells.Item(50, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(50, 3).Value = "Maule"
$ws.Cells.Item(50, 4).Value = 44574
$ws.Cells.Item(50, 4).Style = $ws.Cells.Item(51, 4).Style
$ws.Cells.Item(50, 4).NumberFormat = $ws.Cells.Item(51, 4).NumberFormat
$ws.Cells.Item(50, 5).Value = 7
$ws.Cells.Item(50, 6).Value = 100112001
$ws.Cells.Item(50, 7).Value = "Berenjena"
$ws.Cells.Item(50, 8).Value = "Sin especificar"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 150
$ws.Cells.Item(50, 11).Value = 8000
$ws.Cells.Item(50, 12).Value = 8000
$ws.Cells.Item(50, 13).Value = 8000
$ws.Cells.Item(50, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(50, 15).Value = "Región del Maule"
$ws.Cells.Item(50, 16).Value = 160
$ws.Cells.Item(50, 17).Value = 50
$ws.Cells.Item(50, 18).Value = "Hortaliza"
